# Auto-generated edit script: applies the stock-report correction diff.
# Each statement sets a single cell to its final (post-edit) value, matching
# the literal <v> nodes from the target OOXML (quantities corrected, their
# "Value" (rate*qty) recomputed, Sub Total / Grand Total rows re-summed, and
# a handful of duplicate-code rows with swapped data restored to the right order).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 14
$ws.Range("G5").Value = 5275.34
$ws.Range("F6").Value = 54
$ws.Range("G6").Value = 1613.52
$ws.Range("B10").Value = 26451.15
$ws.Range("F64").Value = 113
$ws.Range("G64").Value = 9172.209999999999
$ws.Range("F71").Value = 303
$ws.Range("G71").Value = 19301.1
$ws.Range("F78").Value = 38
$ws.Range("G78").Value = 2162.2
$ws.Range("F83").Value = 102
$ws.Range("G83").Value = 15368.34
$ws.Range("F84").Value = 25
$ws.Range("G84").Value = 2561.5
$ws.Range("F86").Value = 48
$ws.Range("G86").Value = 6022.56
$ws.Range("B90").Value = 168106.55
$ws.Range("F115").Value = 183
$ws.Range("G115").Value = 17716.23
$ws.Range("B117").Value = 11557.29
$ws.Range("F144").Value = 960
$ws.Range("G144").Value = 8112
$ws.Range("F145").Value = 386
$ws.Range("G145").Value = 3084.14
$ws.Range("B147").Value = 12627.37
$ws.Range("F149").Value = 216
$ws.Range("G149").Value = 13996.8
$ws.Range("F150").Value = 28
$ws.Range("G150").Value = 1301.72
$ws.Range("F151").Value = 89
$ws.Range("G151").Value = 7732.32
$ws.Range("F152").Value = 58
$ws.Range("G152").Value = 5120.82
$ws.Range("B156").Value = 29291.53
$ws.Range("F160").Value = 12
$ws.Range("G160").Value = 1150.92
$ws.Range("B161").Value = 1190.67
$ws.Range("F187").Value = 20
$ws.Range("G187").Value = 999.4
$ws.Range("B192").Value = 64973
$ws.Range("E192").Value = 35.4
$ws.Range("F192").Value = 2
$ws.Range("G192").Value = 66.59999999999999
$ws.Range("B193").Value = 48706
$ws.Range("E193").Value = 39.8
$ws.Range("F193").Value = -144
$ws.Range("G193").Value = -4795.2
$ws.Range("F203").Value = 52
$ws.Range("G203").Value = 1048.32
$ws.Range("F211").Value = 49
$ws.Range("G211").Value = 4958.8
$ws.Range("F213").Value = 7
$ws.Range("G213").Value = 599.76
$ws.Range("F214").Value = 39
$ws.Range("G214").Value = 3420.3
$ws.Range("B216").Value = 34395
$ws.Range("F218").Value = 3
$ws.Range("G218").Value = 648.66
$ws.Range("F222").Value = 10
$ws.Range("G222").Value = 1449.3
$ws.Range("F223").Value = 11
$ws.Range("G223").Value = 1457.28
$ws.Range("F225").Value = 73
$ws.Range("G225").Value = 8338.790000000001
$ws.Range("F233").Value = 113
$ws.Range("G233").Value = 5383.32
$ws.Range("F249").Value = 136
$ws.Range("G249").Value = 18743.52
$ws.Range("F255").Value = 516
$ws.Range("G255").Value = 88406.28
$ws.Range("F256").Value = 259
$ws.Range("G256").Value = 39153.03
$ws.Range("B260").Value = 169980.51
$ws.Range("F283").Value = 37
$ws.Range("G283").Value = 12634.39
$ws.Range("F288").Value = 36
$ws.Range("G288").Value = 3347.64
$ws.Range("F291").Value = 106
$ws.Range("G291").Value = 4559.06
$ws.Range("F292").Value = 41
$ws.Range("G292").Value = 3414.07
$ws.Range("F293").Value = 28
$ws.Range("G293").Value = 1968.96
$ws.Range("F294").Value = 25
$ws.Range("G294").Value = 1784
$ws.Range("F296").Value = 35
$ws.Range("G296").Value = 742
$ws.Range("F302").Value = 34
$ws.Range("G302").Value = 7170.26
$ws.Range("F303").Value = 22
$ws.Range("G303").Value = 4639.58
$ws.Range("B304").Value = 164778.58
$ws.Range("F306").Value = 48
$ws.Range("G306").Value = 1013.28
$ws.Range("B309").Value = 1435.05
$ws.Range("F313").Value = 13
$ws.Range("G313").Value = 1789.71
$ws.Range("B315").Value = 1797.8
$ws.Range("F320").Value = 39
$ws.Range("G320").Value = 2677.35
$ws.Range("F324").Value = 12
$ws.Range("G324").Value = 3177.96
$ws.Range("F327").Value = 11
$ws.Range("G327").Value = 2766.5
$ws.Range("F328").Value = 36
$ws.Range("G328").Value = 1339.56
$ws.Range("F329").Value = 25
$ws.Range("G329").Value = 4161.75
$ws.Range("B330").Value = 25602.3
$ws.Range("F334").Value = 189
$ws.Range("G334").Value = 9793.98
$ws.Range("F338").Value = 73
$ws.Range("G338").Value = 1730.1
$ws.Range("F342").Value = 137
$ws.Range("G342").Value = 4338.79
$ws.Range("F343").Value = 31
$ws.Range("G343").Value = 2231.07
$ws.Range("F345").Value = 31
$ws.Range("G345").Value = 1903.71
$ws.Range("B346").Value = 23764.53
$ws.Range("B364").Value = 53602
$ws.Range("E364").Value = 15.69
$ws.Range("F364").Value = -231
$ws.Range("G364").Value = -3037.65
$ws.Range("B365").Value = 65068
$ws.Range("E365").Value = 13.97
$ws.Range("F365").Value = 63
$ws.Range("G365").Value = 828.45
$ws.Range("B372").Value = 45706
$ws.Range("E372").Value = 23.58
$ws.Range("F372").Value = -202
$ws.Range("G372").Value = -3985.46
$ws.Range("B373").Value = 64922
$ws.Range("E373").Value = 20.98
$ws.Range("F373").Value = 67
$ws.Range("G373").Value = 1321.91
$ws.Range("B380").Value = 64925
$ws.Range("E380").Value = 13.97
$ws.Range("F380").Value = 111
$ws.Range("G380").Value = 1459.65
$ws.Range("B381").Value = 45709
$ws.Range("E381").Value = 15.69
$ws.Range("F381").Value = -300
$ws.Range("G381").Value = -3945
$ws.Range("F450").Value = 8
$ws.Range("G450").Value = 1109.92
$ws.Range("F454").Value = 47
$ws.Range("G454").Value = 1605.05
$ws.Range("F455").Value = 41
$ws.Range("G455").Value = 2608.01
$ws.Range("B460").Value = 12143.27
$ws.Range("B473").Value = 64830
$ws.Range("E473").Value = 34.9
$ws.Range("F473").Value = 105
$ws.Range("G473").Value = 3447.15
$ws.Range("B474").Value = 60022
$ws.Range("E474").Value = 37.22
$ws.Range("F474").Value = -113
$ws.Range("G474").Value = -3709.79
$ws.Range("B475").Value = 45156.9
$ws.Range("F486").Value = 71
$ws.Range("G486").Value = 6272.14
$ws.Range("B488").Value = 28632.56
$ws.Range("F509").Value = 193
$ws.Range("G509").Value = 15513.34
$ws.Range("B510").Value = 20918.22
$ws.Range("F552").Value = 11
$ws.Range("G552").Value = 1119.69
$ws.Range("F555").Value = 15
$ws.Range("G555").Value = 1043.4
$ws.Range("B560").Value = 3648.75
$ws.Range("F577").Value = 40
$ws.Range("G577").Value = 1719.6
$ws.Range("F578").Value = 58
$ws.Range("G578").Value = 2893.62
$ws.Range("B583").Value = 12911.04
$ws.Range("F599").Value = 1325
$ws.Range("G599").Value = 216120.75
$ws.Range("F601").Value = 361
$ws.Range("G601").Value = 102116.07
$ws.Range("F602").Value = 314
$ws.Range("G602").Value = 45420.1
$ws.Range("B606").Value = 364504.97
$ws.Range("F613").Value = 127
$ws.Range("G613").Value = 20213.32
$ws.Range("B618").Value = 41471.83
$ws.Range("B619").Value = 1580364.02
$ws.Range("B620").Value = 1580364.02
